$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Apio" data block (rows 1225-1226),
# pushing the existing rows 1225:1289 down to 1227:1291.
$ws.Range("A1225:A1226").EntireRow.Insert()

# New row 1225: Región Metropolitana, "Primera" quality, week of 2023-12-07
$ws.Cells.Item(1225,1).Value = 6
$ws.Cells.Item(1225,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1225,3).Value = "Metropolitana"
$ws.Cells.Item(1225,4).Value = 45267
$ws.Cells.Item(1225,5).Value = 13
$ws.Cells.Item(1225,6).Value = 100112017
$ws.Cells.Item(1225,7).Value = "Apio"
$ws.Cells.Item(1225,8).Value = "Americana (o)"
$ws.Cells.Item(1225,9).Value = "Primera"
$ws.Cells.Item(1225,10).Value = 1800
$ws.Cells.Item(1225,11).Value = 8000
$ws.Cells.Item(1225,12).Value = 9000
$ws.Cells.Item(1225,13).Value = 8556
$ws.Cells.Item(1225,14).Value = "`$/docena de matas"
$ws.Cells.Item(1225,15).Value = "Región Metropolitana"
$ws.Cells.Item(1225,16).Value = 1426
$ws.Cells.Item(1225,17).Value = 6
$ws.Cells.Item(1225,18).Value = "Hortaliza"

# New row 1226: Región Metropolitana, "Segunda" quality, week of 2023-12-07
$ws.Cells.Item(1226,1).Value = 6
$ws.Cells.Item(1226,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1226,3).Value = "Metropolitana"
$ws.Cells.Item(1226,4).Value = 45267
$ws.Cells.Item(1226,5).Value = 13
$ws.Cells.Item(1226,6).Value = 100112017
$ws.Cells.Item(1226,7).Value = "Apio"
$ws.Cells.Item(1226,8).Value = "Americana (o)"
$ws.Cells.Item(1226,9).Value = "Segunda"
$ws.Cells.Item(1226,10).Value = 900
$ws.Cells.Item(1226,11).Value = 6000
$ws.Cells.Item(1226,12).Value = 7000
$ws.Cells.Item(1226,13).Value = 6333
$ws.Cells.Item(1226,14).Value = "`$/docena de matas"
$ws.Cells.Item(1226,15).Value = "Región Metropolitana"
$ws.Cells.Item(1226,16).Value = 1056
$ws.Cells.Item(1226,17).Value = 6
$ws.Cells.Item(1226,18).Value = "Hortaliza"
